$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at the top of the data (row 2), pushing existing rows down.
$ws.Rows.Item(2).Resize(4).Insert()

# New data for the most recent 4 days, newest first.
$newData = @(
    @("11/23/2023", "15.3% Fat", "70.7 kg"),
    @("11/22/2023", "15.4% Fat", "71 kg"),
    @("11/21/2023", "15.6% Fat", "71.4 kg"),
    @("11/20/2023", "15.7% Fat", "71.6 kg")
)

for ($i = 0; $i -lt $newData.Length; $i++) {
    $row = 2 + $i
    # Force the date column to be treated as text (not auto-converted to a
    # date serial number), then restore the default "Normal" style so the
    # cell's formatting matches the rest of the sheet.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $newData[$i][0]
    $ws.Cells.Item($row, 1).Style = "Normal"
    $ws.Cells.Item($row, 2).Value = $newData[$i][1]
    $ws.Cells.Item($row, 3).Value = $newData[$i][2]
}

# Remove the oldest 4 rows so the table stays the same overall size (41 rows).
$ws.Rows.Item(42).Resize(4).Delete()
